$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.53256663092831
$ws.Range("C2").Value = 12.40119491845101
$ws.Range("E2").Value = 12.15952497240162
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 21.50471880673152
$ws.Range("H2").Value = 12.61806794061573
$ws.Range("L2").Value = 9.584698117903272
$ws.Range("O2").Value = 18.21342026214248
$ws.Range("B3").Value = 15.79218642281184
$ws.Range("C3").Value = 12.23470669892268
$ws.Range("E3").Value = 12.23097805617384
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 21.73691317286601
$ws.Range("H3").Value = 12.69499044762454
$ws.Range("L3").Value = 9.545289273672935
$ws.Range("O3").Value = 18.35892076913278
$ws.Range("B4").Value = 15.31907331714712
$ws.Range("C4").Value = 12.13249652536483
$ws.Range("E4").Value = 12.27763739334332
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 21.89135161574668
$ws.Range("H4").Value = 12.74499673027374
$ws.Range("L4").Value = 9.522473272085916
$ws.Range("O4").Value = 18.45401640277892
$ws.Range("B5").Value = 15.12184286429239
$ws.Range("C5").Value = 12.09089186217784
$ws.Range("E5").Value = 12.29735240715722
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 21.95724575529599
$ws.Range("H5").Value = 12.76607308888132
$ws.Range("L5").Value = 9.513529384440588
$ws.Range("O5").Value = 18.49421398801358
$ws.Range("B6").Value = 15.08883232846119
$ws.Range("C6").Value = 12.08398753459363
$ws.Range("E6").Value = 12.30066841200212
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 21.96836541793327
$ws.Range("H6").Value = 12.76961499746216
$ws.Range("L6").Value = 9.512065810998498
$ws.Range("O6").Value = 18.5009759915482
$ws.Range("B7").Value = 15.31643102829101
$ws.Range("C7").Value = 12.13193518351069
$ws.Range("E7").Value = 12.27790043851238
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 21.89222833887992
$ws.Range("H7").Value = 12.74527814502867
$ws.Range("L7").Value = 9.52235121078261
$ws.Range("O7").Value = 18.45455267253423
$ws.Range("B8").Value = 16.28125057324773
$ws.Range("C8").Value = 12.34381265798646
$ws.Range("E8").Value = 12.18358368677508
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 21.5823002010233
$ws.Range("H8").Value = 12.64401502608481
$ws.Range("L8").Value = 9.570827299818637
$ws.Range("O8").Value = 18.26239229292524
$ws.Range("B9").Value = 18.01810676989871
$ws.Range("C9").Value = 12.7575387190361
$ws.Range("E9").Value = 12.02073326009176
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 21.0698871237667
$ws.Range("H9").Value = 12.46744025861367
$ws.Range("L9").Value = 9.676545847678506
$ws.Range("O9").Value = 17.93136621087561
$ws.Range("B10").Value = 19.19049582108378
$ws.Range("C10").Value = 13.0579580367112
$ws.Range("E10").Value = 11.9145515811695
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 20.75318691295129
$ws.Range("H10").Value = 12.35109469276015
$ws.Range("L10").Value = 9.76028742443853
$ws.Range("O10").Value = 17.71624544440904
$ws.Range("B11").Value = 19.69988415950561
$ws.Range("C11").Value = 13.19335157868471
$ws.Range("E11").Value = 11.86916892334246
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 20.62245359849346
$ws.Range("H11").Value = 12.30106691162463
$ws.Range("L11").Value = 9.799604119399627
$ws.Range("O11").Value = 17.62451826071283
$ws.Range("B12").Value = 19.88923890326313
$ws.Range("C12").Value = 13.24439852879477
$ws.Range("E12").Value = 11.85240362106967
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 20.57489679741344
$ws.Range("H12").Value = 12.28253924752368
$ws.Range("L12").Value = 9.814659379873939
$ws.Range("O12").Value = 17.59066907304032
$ws.Range("B13").Value = 19.84861677134857
$ws.Range("C13").Value = 13.23341524192802
$ws.Range("E13").Value = 11.85599563935397
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 20.58505181360481
$ws.Range("H13").Value = 12.28651097520416
$ws.Range("L13").Value = 9.811409681973178
$ws.Range("O13").Value = 17.59791963316367
$ws.Range("B14").Value = 19.71553392591729
$ws.Range("C14").Value = 13.19755595500465
$ws.Range("E14").Value = 11.8677812103457
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 20.61850182755039
$ws.Range("H14").Value = 12.29953427674994
$ws.Range("L14").Value = 9.800839431893438
$ws.Range("O14").Value = 17.6217156886885
$ws.Range("B15").Value = 19.63355327566645
$ws.Range("C15").Value = 13.17556078302307
$ws.Range("E15").Value = 11.87505493357738
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 20.63924569588848
$ws.Range("H15").Value = 12.30756570125212
$ws.Range("L15").Value = 9.794386313790856
$ws.Range("O15").Value = 17.63640695554189
$ws.Range("B16").Value = 19.15671570469214
$ws.Range("C16").Value = 13.04908066526553
$ws.Range("E16").Value = 11.91757622076586
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 20.76200166056627
$ws.Range("H16").Value = 12.35442244666215
$ws.Range("L16").Value = 9.757741827579286
$ws.Range("O16").Value = 17.72236377436005
$ws.Range("B17").Value = 18.85798565753657
$ws.Range("C17").Value = 12.97113449171056
$ws.Range("E17").Value = 11.94440967522664
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 20.84074617233418
$ws.Range("H17").Value = 12.38390986524746
$ws.Range("L17").Value = 9.735568432657265
$ws.Range("O17").Value = 17.77666897888191
$ws.Range("B18").Value = 18.68391560902888
$ws.Range("C18").Value = 12.92618487541015
$ws.Range("E18").Value = 11.96011838798421
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 20.88729118859823
$ws.Range("H18").Value = 12.40114305744647
$ws.Range("L18").Value = 9.72293047270116
$ws.Range("O18").Value = 17.80848084368983
$ws.Range("B19").Value = 18.62459555744371
$ws.Range("C19").Value = 12.91094688514121
$ws.Range("E19").Value = 11.96548428197736
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 20.90326482830474
$ws.Range("H19").Value = 12.40702477672483
$ws.Range("L19").Value = 9.71867159557322
$ws.Range("O19").Value = 17.81935076463729
$ws.Range("B20").Value = 18.89001944317548
$ws.Range("C20").Value = 12.97944441110761
$ws.Range("E20").Value = 11.94152476377173
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 20.83223375087525
$ws.Range("H20").Value = 12.38074264475911
$ws.Range("L20").Value = 9.737916928719159
$ws.Range("O20").Value = 17.77082835628432
$ws.Range("B21").Value = 19.75472034117375
$ws.Range("C21").Value = 13.20809508762765
$ws.Range("E21").Value = 11.86430809821995
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 20.60862359713687
$ws.Range("H21").Value = 12.29569770622865
$ws.Range("L21").Value = 9.803939713399695
$ws.Range("O21").Value = 17.61470213357643
$ws.Range("B22").Value = 20.29918703101402
$ws.Range("C22").Value = 13.35620989007703
$ws.Range("E22").Value = 11.81629149923928
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 20.47385703626773
$ws.Range("H22").Value = 12.2425452109826
$ws.Range("L22").Value = 9.848057848776334
$ws.Range("O22").Value = 17.51783039848803
$ws.Range("B23").Value = 20.01051382292396
$ws.Range("C23").Value = 13.27729239072123
$ws.Range("E23").Value = 11.84169469033798
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 20.54473309217413
$ws.Range("H23").Value = 12.27069141001889
$ws.Range("L23").Value = 9.824425552117159
$ws.Range("O23").Value = 17.56905859614535
$ws.Range("B24").Value = 18.8755442066299
$ws.Range("C24").Value = 12.97568792113023
$ws.Range("E24").Value = 11.94282815446562
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 20.83607825134334
$ws.Range("H24").Value = 12.38217367143953
$ws.Range("L24").Value = 9.736854830905951
$ws.Range("O24").Value = 17.77346706086871
$ws.Range("B25").Value = 17.56588669230903
$ws.Range("C25").Value = 12.64606206455212
$ws.Range("E25").Value = 12.06242346371376
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 21.19813309852075
$ws.Range("H25").Value = 12.51285612967776
$ws.Range("L25").Value = 9.6468479601913
$ws.Range("O25").Value = 18.01599815863602
